$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old D:T year columns to E:U)
$ws.Columns("D:D").Insert()

# Header for the new column
$ws.Range("D1").Value = "tag"

# Tag every data row with "input-cost"
$ws.Range("D2:D33").Value = "input-cost"

# Give the new column a width similar to column C (closest achievable value;
# the interop engine quantizes stored widths to steps of 1/6 character unit)
$ws.Columns("D:D").ColumnWidth = 15.35

# Update the view: zoom out and select the new tag column's data range
$excel.ActiveWindow.Zoom = 40
$ws.Range("D2:D33").Select() | Out-Null
